$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 24,5
$bf[0,0] = 1.019999999999999
$bf[0,1] = 1.037519885721228
$bf[0,2] = 1.044424457115568
$bf[0,3] = 1.050782331667564
$bf[0,4] = 1.057584328885758
$bf[1,0] = 1.02
$bf[1,1] = 1.038410218801321
$bf[1,2] = 1.04508886747814
$bf[1,3] = 1.051614634263806
$bf[1,4] = 1.058400274699204
$bf[2,0] = 1.02
$bf[2,1] = 1.038987080238854
$bf[2,2] = 1.045519262808981
$bf[2,3] = 1.05215422611793
$bf[2,4] = 1.058929134409287
$bf[3,0] = 1.02
$bf[3,1] = 1.039229772360095
$bf[3,2] = 1.045700313531452
$bf[3,3] = 1.052381316824188
$bf[3,4] = 1.059151677777257
$bf[4,0] = 1.02
$bf[4,1] = 1.039270531921559
$bf[4,2] = 1.045730719291616
$bf[4,3] = 1.052419460742418
$bf[4,4] = 1.059189056091174
$bf[5,0] = 1.02
$bf[5,1] = 1.038990322399142
$bf[5,2] = 1.045521681577668
$bf[5,3] = 1.052157259549605
$bf[5,4] = 1.058932107218952
$bf[6,0] = 1.02
$bf[6,1] = 1.037820620427025
$bf[6,2] = 1.044648897799303
$bf[6,3] = 1.051063396528754
$bf[6,4] = 1.057859896546268
$bf[7,0] = 1.02
$bf[7,1] = 1.035765316134513
$bf[7,2] = 1.043114675559138
$bf[7,3] = 1.049143888809025
$bf[7,4] = 1.055977410667095
$bf[8,0] = 1.02
$bf[8,1] = 1.03439914307695
$bf[8,2] = 1.042094476668717
$bf[8,3] = 1.047869717421468
$bf[8,4] = 1.054727162926161
$bf[9,0] = 1.02
$bf[9,1] = 1.03380854970372
$bf[9,2] = 1.041653361968386
$bf[9,3] = 1.047319314788292
$bf[9,4] = 1.054186940610297
$bf[10,0] = 1.02
$bf[10,1] = 1.033589324048416
$bf[10,2] = 1.041489609999868
$bf[10,3] = 1.0471150712191
$bf[10,4] = 1.053986451669231
$bf[11,0] = 1.02
$bf[11,1] = 1.033636342018542
$bf[11,2] = 1.0415247309043
$bf[11,3] = 1.047158873056181
$bf[11,4] = 1.054029449337331
$bf[12,0] = 1.02
$bf[12,1] = 1.033790425434804
$bf[12,2] = 1.041639824176726
$bf[12,3] = 1.047302427850743
$bf[12,4] = 1.054170364577717
$bf[13,0] = 1.02
$bf[13,1] = 1.033885380839505
$bf[13,2] = 1.041710749936597
$bf[13,3] = 1.047390903318214
$bf[13,4] = 1.054257210180071
$bf[14,0] = 1.02
$bf[14,1] = 1.034438359348951
$bf[14,2] = 1.042123765597138
$bf[14,3] = 1.047906273843191
$bf[14,4] = 1.054763039952597
$bf[15,0] = 1.02
$bf[15,1] = 1.034785488532662
$bf[15,2] = 1.042383011529582
$bf[15,3] = 1.048229907405255
$bf[15,4] = 1.055080641034093
$bf[16,0] = 1.02
$bf[16,1] = 1.034988056465123
$bf[16,2] = 1.042534286747863
$bf[16,3] = 1.048418804814166
$bf[16,4] = 1.055266002545205
$bf[17,0] = 1.02
$bf[17,1] = 1.03505714268385
$bf[17,2] = 1.042585878056623
$bf[17,3] = 1.048483235506339
$bf[17,4] = 1.055329224659065
$bf[18,0] = 1.02
$bf[18,1] = 1.03474823518971
$bf[18,2] = 1.042355190525277
$bf[18,3] = 1.048195171393013
$bf[18,4] = 1.055046554034115
$bf[19,0] = 1.02
$bf[19,1] = 1.033745047638595
$bf[19,2] = 1.041605929365983
$bf[19,3] = 1.047260148992698
$bf[19,4] = 1.054128863741744
$bf[20,0] = 1.02
$bf[20,1] = 1.033115155166347
$bf[20,2] = 1.041135405123982
$bf[20,3] = 1.046673424109418
$bf[20,4] = 1.053552881246684
$bf[21,0] = 1.02
$bf[21,1] = 1.033448991882183
$bf[21,2] = 1.041384784678829
$bf[21,3] = 1.04698434741808
$bf[21,4] = 1.053858124482468
$bf[22,0] = 1.02
$bf[22,1] = 1.034765068083564
$bf[22,2] = 1.042367761449933
$bf[22,3] = 1.048210866708253
$bf[22,4] = 1.05506195614197
$bf[23,0] = 1.02
$bf[23,1] = 1.036295957874036
$bf[23,2] = 1.04351085540309
$bf[23,3] = 1.049639166066543
$bf[23,4] = 1.056463250693674

$ws.Range("B2:F25").Value2 = $bf

$inn = New-Object 'object[,]' 24,6
$inn[0,0] = 1.037148687827609
$inn[0,1] = 1.04262231252285
$inn[0,2] = 1.047195545589036
$inn[0,3] = 1.053535653057078
$inn[0,4] = 1.060318898330117
$inn[0,5] = 1.018125825952509
$inn[1,0] = 1.037282031760202
$inn[1,1] = 1.043157355805672
$inn[1,2] = 1.047671574350511
$inn[1,3] = 1.054180419920195
$inn[1,4] = 1.060948703901473
$inn[1,5] = 1.018305557543465
$inn[2,0] = 1.037367190653197
$inn[2,1] = 1.043503674883248
$inn[2,2] = 1.047979390252241
$inn[2,3] = 1.054598013342429
$inn[2,4] = 1.061356463884899
$inn[2,5] = 1.018421833660802
$inn[3,0] = 1.037402721813344
$inn[3,1] = 1.043649292559623
$inn[3,2] = 1.048108745298317
$inn[3,3] = 1.054773660551483
$inn[3,4] = 1.061527940557921
$inn[3,5] = 1.018470710334918
$inn[4,0] = 1.03740867182062
$inn[4,1] = 1.04367374384535
$inn[4,2] = 1.048130461565026
$inn[4,3] = 1.054803157779305
$inn[4,4] = 1.061556735389025
$inn[4,5] = 1.018478916589131
$inn[5,0] = 1.037367666482614
$inn[5,1] = 1.04350562053659
$inn[5,2] = 1.047981118901794
$inn[5,3] = 1.054600359993618
$inn[5,4] = 1.0613587549529
$inn[5,5] = 1.018422486776729
$inn[6,0] = 1.037193984275485
$inn[6,1] = 1.042803109304641
$inn[6,2] = 1.047356463713164
$inn[6,3] = 1.053753474174816
$inn[6,4] = 1.060531694749562
$inn[6,5] = 1.018186571400911
$inn[7,0] = 1.036879358170115
$inn[7,1] = 1.041566095616009
$inn[7,2] = 1.046254217695499
$inn[7,3] = 1.052264177371582
$inn[7,4] = 1.059076167749025
$inn[7,5] = 1.017770709299369
$inn[8,0] = 1.036663880825162
$inn[8,1] = 1.040742101087567
$inn[8,2] = 1.045518440787194
$inn[8,3] = 1.051273436133985
$inn[8,4] = 1.058107159344709
$inn[8,5] = 1.017493396324403
$inn[9,0] = 1.036569226753983
$inn[9,1] = 1.040385479855142
$inn[9,2] = 1.045199633322448
$inn[9,3] = 1.050844956820723
$inn[9,4] = 1.057687906183367
$inn[9,5] = 1.017373305803039
$inn[10,0] = 1.036533865700757
$inn[10,1] = 1.040253042243656
$inn[10,2] = 1.045081183625847
$inn[10,3] = 1.050685879902941
$inn[10,4] = 1.057532228633958
$inn[10,5] = 1.01732869749181
$inn[11,0] = 1.036541459911524
$inn[11,1] = 1.040281449300816
$inn[11,2] = 1.045106592842708
$inn[11,3] = 1.050719998825307
$inn[11,4] = 1.057565619629942
$inn[11,5] = 1.017338266179066
$inn[12,0] = 1.036566307921731
$inn[12,1] = 1.040374531954676
$inn[12,2] = 1.045189842846269
$inn[12,3] = 1.050831805844975
$inn[12,4] = 1.057675036763903
$inn[12,5] = 1.017369618490746
$inn[13,0] = 1.036581590811433
$inn[13,1] = 1.040431886917823
$inn[13,2] = 1.045241131931086
$inn[13,3] = 1.050900704403804
$inn[13,4] = 1.057742459174974
$inn[13,5] = 1.017388935521044
$inn[14,0] = 1.036670134296395
$inn[14,1] = 1.040765772635679
$inn[14,2] = 1.045539594651625
$inn[14,3] = 1.051301883942601
$inn[14,4] = 1.058134990965544
$inn[14,5] = 1.017501366121669
$inn[15,0] = 1.036725314053291
$inn[15,1] = 1.040975257766685
$inn[15,2] = 1.045726756818037
$inn[15,3] = 1.051553673077388
$inn[15,4] = 1.058381306290674
$inn[15,5] = 1.017571887949932
$inn[16,0] = 1.03675736908034
$inn[16,1] = 1.041097463586439
$inn[16,2] = 1.045835904875977
$inn[16,3] = 1.051700587314297
$inn[16,4] = 1.058525009960265
$inn[16,5] = 1.01761302093222
$inn[17,0] = 1.036768276879771
$inn[17,1] = 1.041139135404497
$inn[17,2] = 1.045873118031275
$inn[17,3] = 1.05175068968902
$inn[17,4] = 1.058574014595883
$inn[17,5] = 1.017627045996741
$inn[18,0] = 1.03671940726992
$inn[18,1] = 1.040952780259447
$inn[18,2] = 1.045706678196341
$inn[18,3] = 1.051526653314294
$inn[18,4] = 1.058354875654923
$inn[18,5] = 1.017564321746522
$inn[19,0] = 1.036558996381754
$inn[19,1] = 1.040347120666714
$inn[19,2] = 1.045165328623999
$inn[19,3] = 1.050798879260941
$inn[19,4] = 1.057642814704071
$inn[19,5] = 1.017360386049621
$inn[20,0] = 1.036456969470095
$inn[20,1] = 1.039966477560529
$inn[20,2] = 1.044824785377517
$inn[20,3] = 1.050341758836162
$inn[20,4] = 1.05719541386851
$inn[20,5] = 1.017232156043117
$inn[21,0] = 1.036511166572402
$inn[21,1] = 1.040168248164336
$inn[21,2] = 1.045005329969158
$inn[21,3] = 1.050584042909142
$inn[21,4] = 1.05743256053812
$inn[21,5] = 1.017300133748777
$inn[22,0] = 1.036722076694037
$inn[22,1] = 1.040962936826028
$inn[22,2] = 1.045715750924612
$inn[22,3] = 1.051538862228205
$inn[22,4] = 1.058366818422511
$inn[22,5] = 1.017567740591949
$inn[23,0] = 1.037148687827609
$inn[23,1] = 1.041885778763931
$inn[23,2] = 1.046539347221533
$inn[23,3] = 1.052648828637528
$inn[23,4] = 1.059452226491994
$inn[23,5] = 1.017878234332989

$ws.Range("I2:N25").Value2 = $inn

Write-Host "Updated vm_pu data for rows 2-25 (380 kV case)"
